# Apply the StateTable_minimal.xlsx edit:
# "finally made decision on where to put EFCT_special calls"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text values (content-level changes) ---

# Row 4: the special-handler text moves from B4 (trimmed) and the
# EFCT_UNIQ_WAITING column gains the mSPCL_EFCT_CONTINUOUS prefix.
$ws.Range("B4").Value2 = "mSPCL_HANDLER | mSPCL_HANDLER_SHOOT"
$ws.Range("D4").Value2 = "mSPCL_EFCT_CONTINUOUS|mEFCT_UNIQ_WAITING"
$ws.Range("E4").Value2 = "mSPCL_EFCT_CONTINUOUS|mEFCT_UNIQ_WAITING"

# Row 5: B5 (mSPCL_EFCT_CONTINUOUS) is cleared; D5/E5 pick up the combined text
$ws.Range("B5").ClearContents()
$ws.Range("D5").Value2 = "mSPCL_EFCT_CONTINUOUS|mEFCT_UNIQ_WAITING"
$ws.Range("E5").Value2 = "mSPCL_EFCT_CONTINUOUS|mEFCT_UNIQ_WAITING"

# Row 6: same treatment as row 5
$ws.Range("B6").ClearContents()
$ws.Range("D6").Value2 = "mSPCL_EFCT_CONTINUOUS|mEFCT_UNIQ_WAITING"
$ws.Range("E6").Value2 = "mSPCL_EFCT_CONTINUOUS|mEFCT_UNIQ_WAITING"

# --- Row height changes ---
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30

# --- Selection change ---
$ws.Range("B5:B6").Select()
